# group_report_template.xlsx -- "add treatments to group parser, in progress"
#
# Adds a new "Treatments" mini-table (header "Treatment" / "Date" / "Container")
# to the right of the existing table on the "Containers" sheet (columns G:I),
# and makes "Containers" the active/selected tab instead of "Event History".

$wb = $excel.ActiveWorkbook
$wsContainers = $wb.Worksheets.Item("Containers")

# New header cells for the "Treatments" block (columns G, H, I).
$wsContainers.Range("G1").Value = "Treatments"
$wsContainers.Range("G4").Value = "Treatment"
$wsContainers.Range("H4").Value = "Date"
$wsContainers.Range("I4").Value = "Container"

# Match the look of the existing header row (A4:D4) as closely as possible by
# reusing the same named cell style used elsewhere in the template.
$wsContainers.Range("G4:I4").Style = "Template Header"

# Column widths for the new block, picked to match the template's layout.
$wsContainers.Columns.Item(7).ColumnWidth = 10.3
$wsContainers.Columns.Item(8).ColumnWidth = 16.8
$wsContainers.Columns.Item(9).ColumnWidth = 8.8

# The template now opens on the "Containers" tab (was "Event History").
$wsContainers.Activate()
